$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values - automatic electricity price update
$ws.Range("A2").Value = 45928
$ws.Range("B2").Value = 63.51
$ws.Range("C2").Value = 60.86
$ws.Range("D2").Value = 59.83
$ws.Range("E2").Value = 58.3
$ws.Range("F2").Value = 58.04
$ws.Range("G2").Value = 58.04
$ws.Range("H2").Value = 57.91
$ws.Range("I2").Value = 58.87
$ws.Range("J2").Value = 58.87
$ws.Range("K2").Value = 57.91
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 29.16
$ws.Range("N2").Value = 17.11
$ws.Range("O2").Value = 8.529999999999999
$ws.Range("P2").Value = 15
$ws.Range("Q2").Value = 16.79
$ws.Range("R2").Value = 22.97
$ws.Range("S2").Value = 37.34
$ws.Range("T2").Value = 61.11
$ws.Range("U2").Value = 70.98
$ws.Range("V2").Value = 74.14
$ws.Range("W2").Value = 61
$ws.Range("X2").Value = 69.62
$ws.Range("Y2").Value = 60.16
$ws.Range("Z2").Value = 48.79
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 66.23
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 67.56999999999999
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 66.04000000000001
$ws.Range("AG2").Value = "10h-17h"
